# Update DateBase/orders/International Ever Green_2025-11-7.xlsx
#
# 1) Orders sheet (sheet1): F71 "1" -> "10", and append 4 new rows (72-75)
#    of flower-order data. All of these cells hold numeric-looking text,
#    so they are entered with a leading apostrophe (forcing text, matching
#    the workbook's existing "number stored as text" convention) and the
#    style is then reset to Normal so no stray number-format / quote-prefix
#    style is introduced.
# 2) Summary sheet (sheet2): G2 gets the extra run lengths appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("F71") "10"

Set-TextValue $ws.Range("C72") "11_香槟洋桔梗_Champagne Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
Set-TextValue $ws.Range("F72") "8"

Set-TextValue $ws.Range("A73") "5"
Set-TextValue $ws.Range("C73") "11_香槟洋桔梗_Champagne Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
Set-TextValue $ws.Range("F73") "12"

Set-TextValue $ws.Range("C74") "8_冰淇淋洋桔梗_Icecream Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
Set-TextValue $ws.Range("F74") "10"

Set-TextValue $ws.Range("C75") "14_波浪浅紫洋桔梗_Wavy Light Purple Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
Set-TextValue $ws.Range("F75") "10"

$summary = $wb.Worksheets.Item("Summary")
$summary.Range("G2").Value = "0202026271350151315142075625361010341035201420830208540445595235361550351691515258101068.5101.51410710510201010301010202510201010108121010"
